# Update the "Förändrad" (Changed) date column (C) for rows 2-39
# from 2023-09-21 (45190) to 2023-09-23 (45192).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 39; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
